# Apply updated crypto price/volume values per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "37.387.64"
Set-TextValue "E2" "  -0.09%  "
Set-TextValue "D3" "2.068.98"
Set-TextValue "E4" "  +0.06%  "
Set-TextValue "D5" "235.48"
Set-TextValue "E5" "  -0.23%  "
Set-TextValue "D6" "0.626"
Set-TextValue "E6" "  +1.81%  "
Set-TextValue "E7" "  +0.06%  "
Set-TextValue "D8" "57.30"
Set-TextValue "E8" "  -1.72%  "
Set-TextValue "E9" "  +2.95%  "
Set-TextValue "E10" "  +1.22%  "
Set-TextValue "D11" "0.103"
Set-TextValue "E11" "  +0.82%  "
Set-TextValue "D12" "2.373.44"
Set-TextValue "E12" "  +0.32%  "
Set-TextValue "E13" "  -1.19%  "
Set-TextValue "D14" "20.71"
Set-TextValue "E14" "  -1.59%  "
Set-TextValue "D15" "0.776"
Set-TextValue "E15" "  -0.31%  "
Set-TextValue "E16" "  -0.73%  "
Set-TextValue "D17" "2.067.09"
Set-TextValue "E17" "  +0.18%  "
Set-TextValue "D18" "37.332.06"
Set-TextValue "E18" "  -0.50%  "
Set-TextValue "D19" "6.20"
Set-TextValue "E19" "  -0.46%  "
Set-TextValue "D20" "69.59"
Set-TextValue "E20" "  +0.72%  "
Set-TextValue "E21" "  -0.11%  "
Set-TextValue "D22" "226.31"
Set-TextValue "E22" "  -0.21%  "
Set-TextValue "D24" "2.44"
Set-TextValue "E24" "  +1.88%  "
Set-TextValue "E25" "  -1.93%  "
Set-TextValue "D26" "166.78"
Set-TextValue "E26" "  +1.37%  "
Set-TextValue "D27" "8.95"
Set-TextValue "E27" "  +0.84%  "
Set-TextValue "D28" "1.40"
Set-TextValue "E28" "  -6.26%  "
Set-TextValue "D30" "19.10"
Set-TextValue "E30" "  -0.75%  "
Set-TextValue "D31" "0.118"
Set-TextValue "E31" "  -1.33%  "
Set-TextValue "E32" "  +0.78%  "
Set-TextValue "D33" "0.0617"
Set-TextValue "E33" "  -1.25%  "
Set-TextValue "E34" "  +1.08%  "
Set-TextValue "D36" "1.79"
Set-TextValue "E36" "  +0.46%  "
Set-TextValue "E37" "  -2.40%  "
Set-TextValue "E38" "  +0.00%  "
Set-TextValue "E39" "  -5.00%  "
Set-TextValue "E40" "  -0.97%  "
Set-TextValue "D41" "0.0959"
Set-TextValue "E41" "  -2.64%  "
Set-TextValue "D42" "97.75"
Set-TextValue "E42" "  +0.72%  "
Set-TextValue "D43" "1.478.83"
Set-TextValue "E43" "  +0.34%  "
Set-TextValue "D44" "0.0213"
Set-TextValue "E44" "  +0.84%  "
Set-TextValue "E45" "  -0.43%  "
Set-TextValue "D46" "4.18"
Set-TextValue "E46" "  -6.80%  "
Set-TextValue "D47" "1.03"
Set-TextValue "E47" "  -0.20%  "
Set-TextValue "D48" "7.21"
Set-TextValue "E48" "  -0.95%  "
Set-TextValue "D49" "15.10"
Set-TextValue "E49" "  -5.26%  "
Set-TextValue "E50" "  +0.90%  "
Set-TextValue "D51" "2.259.88"
Set-TextValue "E51" "  +0.26%  "
